# fixup style mistake in Document Principles
#
# The three "Feed updates to ..." text boxes (TextBox 23, TextBox 24,
# TextBox 33 -- shape ids 24, 25, 34) were originally stacked below the
# "Rectangle: Folded Corner 25" shape (and the other shapes that follow
# it) in the slide's z-order/document order. They need to be moved to
# the very end of the shape tree (on top of everything else), and their
# semi-transparent accent1 fill needs to go from 50% opaque (alpha
# 50000 / Fill.Transparency 0.5) to 93% opaque (alpha 93000 /
# Fill.Transparency 0.07).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# In document order (before the edit) these three boxes are shapes
# 17, 18, 19. Repeatedly grabbing whatever is now at position 17 and
# sending it to the front of the z-order (ZOrder msoBringToFront = 0)
# moves TextBox 23, then TextBox 24, then TextBox 33 to the end of the
# shape tree, preserving their relative order, while also bumping the
# fill transparency on each.
for ($i = 0; $i -lt 3; $i++) {
  $shp = $s.Shapes.Item(17)
  $shp.ZOrder(0)
  $shp.Fill.Transparency = 0.07
}
